$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.267.98"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "2.050.65"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.58"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  +2.62%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.60"
$ws.Range("E8").Value = "  +5.11%  "
$ws.Range("E9").Value = "  +3.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.64"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0758"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.101"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "2.362.60"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.32"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.82"
$ws.Range("E15").Value = "  +3.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.771"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.16"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "2.056.02"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("D19").Value = "37.494.04"
$ws.Range("E19").Value = "  +3.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.14"
$ws.Range("E20").Value = "  +14.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.99"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").Value = "0.0₃0809"
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "225.41"
$ws.Range("E23").Value = "  +2.30%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.39"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.42"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.20"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.47"
$ws.Range("E28").Value = "  +6.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.87"
$ws.Range("E29").Value = "  +2.81%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.06"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.49"
$ws.Range("E33").Value = "  +3.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0619"
$ws.Range("E34").Value = "  +2.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.55"
$ws.Range("E35").Value = "  +4.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.58"
$ws.Range("E36").Value = "  +7.90%  "
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.92"
$ws.Range("E38").Value = "  +3.92%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.30"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.76"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.64"
$ws.Range("E41").Value = "  +7.59%  "
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0941"
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("D44").Value = "1.454.84"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("E45").Value = "  +5.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "95.41"
$ws.Range("E46").Value = "  +5.98%  "
$ws.Range("E47").Value = "  +3.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.61"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.02"
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.15"
$ws.Range("E50").Value = "  +4.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.93"
$ws.Range("E51").Value = "  +1.82%  "
